# Apply the "4 mdelo melhores rstds" edit:
# Column A (model names) is reordered/renumbered for rows 2-26, and the
# metric columns B:I are all set to the same (new, "best") set of values
# for every one of those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$modelNames = @(
    "model_10_3_0",
    "model_10_3_22",
    "model_10_3_21",
    "model_10_3_20",
    "model_10_3_19",
    "model_10_3_18",
    "model_10_3_17",
    "model_10_3_16",
    "model_10_3_15",
    "model_10_3_14",
    "model_10_3_13",
    "model_10_3_23",
    "model_10_3_12",
    "model_10_3_10",
    "model_10_3_9",
    "model_10_3_8",
    "model_10_3_7",
    "model_10_3_6",
    "model_10_3_5",
    "model_10_3_4",
    "model_10_3_3",
    "model_10_3_2",
    "model_10_3_1",
    "model_10_3_11",
    "model_10_3_24"
)

$metricValues = @(
    0.6731329884640765,
    0.6356498070262204,
    -0.4299051501425744,
    0.2435373185495373,
    0.3617455065250397,
    0.3348711133003235,
    0.8773518204689026,
    0.590156614780426
)

$startRow = 2
for ($i = 0; $i -lt $modelNames.Length; $i++) {
    $row = $startRow + $i

    $ws.Cells.Item($row, 1).Value = $modelNames[$i]

    for ($c = 0; $c -lt $metricValues.Length; $c++) {
        $col = 2 + $c
        $ws.Cells.Item($row, $col).Value = $metricValues[$c]
    }
}
